$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values (row -> col letter -> value) before overwriting,
# since this edit permutes existing row data across rows (weekly refresh).
$orig = @{}
$orig[2] = @{
    'D' = 44579
    'K' = 'Modesto'
    'L' = 'Primera'
    'M' = 180
    'N' = 13000
    'O' = 14000
    'P' = 13444
    'Q' = '$/caja 18 kilos'
    'R' = 'Región Metropolitana'
    'S' = 747
    'T' = 18
}
$orig[3] = @{
    'D' = 44944
    'K' = 'Modesto'
    'L' = 'Primera'
    'M' = 100
    'N' = 16000
    'O' = 17000
    'P' = 16500
    'Q' = '$/caja 16 kilos empedrada'
    'R' = 'Región de O''Higgins'
    'S' = 1031
    'T' = 16
}
$orig[4] = @{
    'D' = 44944
    'K' = 'Modesto'
    'L' = 'Segunda'
    'M' = 50
    'N' = 14000
    'O' = 14000
    'P' = 14000
    'Q' = '$/caja 16 kilos empedrada'
    'R' = 'Región de O''Higgins'
    'S' = 875
    'T' = 16
}
$orig[5] = @{
    'D' = 44189
    'K' = 'Dina'
    'L' = 'Primera'
    'M' = 200
    'N' = 15000
    'O' = 16000
    'P' = 15500
    'Q' = '$/caja 15 kilos granel'
    'R' = 'Región de O''Higgins'
    'S' = 1033
    'T' = 15
}
$orig[6] = @{
    'D' = 44189
    'K' = 'Dina'
    'L' = 'Segunda'
    'M' = 100
    'N' = 14000
    'O' = 14000
    'P' = 14000
    'Q' = '$/caja 15 kilos granel'
    'R' = 'Región de O''Higgins'
    'S' = 933
    'T' = 15
}
$orig[7] = @{
    'D' = 44901
    'K' = 'Castle Brite'
    'L' = 'Primera'
    'M' = 100
    'N' = 15000
    'O' = 16000
    'P' = 15500
    'Q' = '$/caja 10 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1550
    'T' = 10
}
$orig[8] = @{
    'D' = 44908
    'K' = 'Albaricoque'
    'L' = 'Primera'
    'M' = 100
    'N' = 20000
    'O' = 22000
    'P' = 21000
    'Q' = '$/caja 18 kilos granel'
    'R' = 'Región de O''Higgins'
    'S' = 1167
    'T' = 18
}
$orig[9] = @{
    'D' = 44918
    'K' = 'Dina'
    'L' = 'Primera'
    'M' = 100
    'N' = 17000
    'O' = 18000
    'P' = 17500
    'Q' = '$/caja 18 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 972
    'T' = 18
}
$orig[10] = @{
    'D' = 44545
    'K' = 'Castle Brite'
    'L' = 'Primera'
    'M' = 100
    'N' = 18000
    'O' = 19000
    'P' = 18500
    'Q' = '$/caja 15 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1233
    'T' = 15
}
$orig[11] = @{
    'D' = 44545
    'K' = 'Castle Brite'
    'L' = 'Segunda'
    'M' = 50
    'N' = 17000
    'O' = 17000
    'P' = 17000
    'Q' = '$/caja 15 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1133
    'T' = 15
}
$orig[12] = @{
    'D' = 44875
    'K' = 'Castle Brite'
    'L' = 'Primera'
    'M' = 50
    'N' = 31000
    'O' = 32000
    'P' = 31400
    'Q' = '$/bandeja 10 kilos'
    'R' = 'Provincia de Limarí'
    'S' = 3140
    'T' = 10
}
$orig[13] = @{
    'D' = 44559
    'K' = 'Modesto'
    'L' = 'Primera'
    'M' = 100
    'N' = 19000
    'O' = 20000
    'P' = 19500
    'Q' = '$/caja 18 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1083
    'T' = 18
}
$orig[14] = @{
    'D' = 44559
    'K' = 'Modesto'
    'L' = 'Segunda'
    'M' = 50
    'N' = 18000
    'O' = 18000
    'P' = 18000
    'Q' = '$/caja 18 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1000
    'T' = 18
}
$orig[15] = @{
    'D' = 44187
    'K' = 'Dina'
    'L' = 'Primera'
    'M' = 100
    'N' = 15000
    'O' = 16000
    'P' = 15500
    'Q' = '$/caja 18 kilos'
    'R' = 'Región Metropolitana'
    'S' = 861
    'T' = 18
}
$orig[16] = @{
    'D' = 44938
    'K' = 'Modesto'
    'L' = 'Primera'
    'M' = 270
    'N' = 14000
    'O' = 15000
    'P' = 14556
    'Q' = '$/caja 15 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 970
    'T' = 15
}
$orig[17] = @{
    'D' = 44902
    'K' = 'Castle Brite'
    'L' = 'Primera'
    'M' = 100
    'N' = 15000
    'O' = 16000
    'P' = 15500
    'Q' = '$/caja 10 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1550
    'T' = 10
}
$orig[18] = @{
    'D' = 44902
    'K' = 'Castle Brite'
    'L' = 'Segunda'
    'M' = 50
    'N' = 13000
    'O' = 13000
    'P' = 13000
    'Q' = '$/caja 10 kilos'
    'R' = 'Región de O''Higgins'
    'S' = 1300
    'T' = 10
}
$orig[19] = @{
    'D' = 44159
    'K' = 'Castle Brite'
    'L' = 'Primera'
    'M' = 100
    'N' = 14000
    'O' = 15000
    'P' = 14500
    'Q' = '$/caja 15 kilos'
    'R' = 'Región Metropolitana'
    'S' = 967
    'T' = 15
}

# Apply the permuted values per the diff (destination row <- source row, using original/pre-edit data)
# Row 2 <- original row 13
$ws.Cells.Item(2, 4).Value = $orig[13]['D']
$ws.Cells.Item(2, 11).Value = $orig[13]['K']
$ws.Cells.Item(2, 12).Value = $orig[13]['L']
$ws.Cells.Item(2, 13).Value = $orig[13]['M']
$ws.Cells.Item(2, 14).Value = $orig[13]['N']
$ws.Cells.Item(2, 15).Value = $orig[13]['O']
$ws.Cells.Item(2, 16).Value = $orig[13]['P']
$ws.Cells.Item(2, 17).Value = $orig[13]['Q']
$ws.Cells.Item(2, 18).Value = $orig[13]['R']
$ws.Cells.Item(2, 19).Value = $orig[13]['S']
$ws.Cells.Item(2, 20).Value = $orig[13]['T']

# Row 3 <- original row 14
$ws.Cells.Item(3, 4).Value = $orig[14]['D']
$ws.Cells.Item(3, 11).Value = $orig[14]['K']
$ws.Cells.Item(3, 12).Value = $orig[14]['L']
$ws.Cells.Item(3, 13).Value = $orig[14]['M']
$ws.Cells.Item(3, 14).Value = $orig[14]['N']
$ws.Cells.Item(3, 15).Value = $orig[14]['O']
$ws.Cells.Item(3, 16).Value = $orig[14]['P']
$ws.Cells.Item(3, 17).Value = $orig[14]['Q']
$ws.Cells.Item(3, 18).Value = $orig[14]['R']
$ws.Cells.Item(3, 19).Value = $orig[14]['S']
$ws.Cells.Item(3, 20).Value = $orig[14]['T']

# Row 4 <- original row 8
$ws.Cells.Item(4, 4).Value = $orig[8]['D']
$ws.Cells.Item(4, 11).Value = $orig[8]['K']
$ws.Cells.Item(4, 12).Value = $orig[8]['L']
$ws.Cells.Item(4, 13).Value = $orig[8]['M']
$ws.Cells.Item(4, 14).Value = $orig[8]['N']
$ws.Cells.Item(4, 15).Value = $orig[8]['O']
$ws.Cells.Item(4, 16).Value = $orig[8]['P']
$ws.Cells.Item(4, 17).Value = $orig[8]['Q']
$ws.Cells.Item(4, 18).Value = $orig[8]['R']
$ws.Cells.Item(4, 19).Value = $orig[8]['S']
$ws.Cells.Item(4, 20).Value = $orig[8]['T']

# Row 5 <- original row 19
$ws.Cells.Item(5, 4).Value = $orig[19]['D']
$ws.Cells.Item(5, 11).Value = $orig[19]['K']
$ws.Cells.Item(5, 12).Value = $orig[19]['L']
$ws.Cells.Item(5, 13).Value = $orig[19]['M']
$ws.Cells.Item(5, 14).Value = $orig[19]['N']
$ws.Cells.Item(5, 15).Value = $orig[19]['O']
$ws.Cells.Item(5, 16).Value = $orig[19]['P']
$ws.Cells.Item(5, 17).Value = $orig[19]['Q']
$ws.Cells.Item(5, 18).Value = $orig[19]['R']
$ws.Cells.Item(5, 19).Value = $orig[19]['S']
$ws.Cells.Item(5, 20).Value = $orig[19]['T']

# Row 6 <- original row 17
$ws.Cells.Item(6, 4).Value = $orig[17]['D']
$ws.Cells.Item(6, 11).Value = $orig[17]['K']
$ws.Cells.Item(6, 12).Value = $orig[17]['L']
$ws.Cells.Item(6, 13).Value = $orig[17]['M']
$ws.Cells.Item(6, 14).Value = $orig[17]['N']
$ws.Cells.Item(6, 15).Value = $orig[17]['O']
$ws.Cells.Item(6, 16).Value = $orig[17]['P']
$ws.Cells.Item(6, 17).Value = $orig[17]['Q']
$ws.Cells.Item(6, 18).Value = $orig[17]['R']
$ws.Cells.Item(6, 19).Value = $orig[17]['S']
$ws.Cells.Item(6, 20).Value = $orig[17]['T']

# Row 7 <- original row 18
$ws.Cells.Item(7, 4).Value = $orig[18]['D']
$ws.Cells.Item(7, 11).Value = $orig[18]['K']
$ws.Cells.Item(7, 12).Value = $orig[18]['L']
$ws.Cells.Item(7, 13).Value = $orig[18]['M']
$ws.Cells.Item(7, 14).Value = $orig[18]['N']
$ws.Cells.Item(7, 15).Value = $orig[18]['O']
$ws.Cells.Item(7, 16).Value = $orig[18]['P']
$ws.Cells.Item(7, 17).Value = $orig[18]['Q']
$ws.Cells.Item(7, 18).Value = $orig[18]['R']
$ws.Cells.Item(7, 19).Value = $orig[18]['S']
$ws.Cells.Item(7, 20).Value = $orig[18]['T']

# Row 8 <- original row 12
$ws.Cells.Item(8, 4).Value = $orig[12]['D']
$ws.Cells.Item(8, 11).Value = $orig[12]['K']
$ws.Cells.Item(8, 12).Value = $orig[12]['L']
$ws.Cells.Item(8, 13).Value = $orig[12]['M']
$ws.Cells.Item(8, 14).Value = $orig[12]['N']
$ws.Cells.Item(8, 15).Value = $orig[12]['O']
$ws.Cells.Item(8, 16).Value = $orig[12]['P']
$ws.Cells.Item(8, 17).Value = $orig[12]['Q']
$ws.Cells.Item(8, 18).Value = $orig[12]['R']
$ws.Cells.Item(8, 19).Value = $orig[12]['S']
$ws.Cells.Item(8, 20).Value = $orig[12]['T']

# Row 10 <- original row 15
$ws.Cells.Item(10, 4).Value = $orig[15]['D']
$ws.Cells.Item(10, 11).Value = $orig[15]['K']
$ws.Cells.Item(10, 12).Value = $orig[15]['L']
$ws.Cells.Item(10, 13).Value = $orig[15]['M']
$ws.Cells.Item(10, 14).Value = $orig[15]['N']
$ws.Cells.Item(10, 15).Value = $orig[15]['O']
$ws.Cells.Item(10, 16).Value = $orig[15]['P']
$ws.Cells.Item(10, 17).Value = $orig[15]['Q']
$ws.Cells.Item(10, 18).Value = $orig[15]['R']
$ws.Cells.Item(10, 19).Value = $orig[15]['S']
$ws.Cells.Item(10, 20).Value = $orig[15]['T']

# Row 11 <- original row 3
$ws.Cells.Item(11, 4).Value = $orig[3]['D']
$ws.Cells.Item(11, 11).Value = $orig[3]['K']
$ws.Cells.Item(11, 12).Value = $orig[3]['L']
$ws.Cells.Item(11, 13).Value = $orig[3]['M']
$ws.Cells.Item(11, 14).Value = $orig[3]['N']
$ws.Cells.Item(11, 15).Value = $orig[3]['O']
$ws.Cells.Item(11, 16).Value = $orig[3]['P']
$ws.Cells.Item(11, 17).Value = $orig[3]['Q']
$ws.Cells.Item(11, 18).Value = $orig[3]['R']
$ws.Cells.Item(11, 19).Value = $orig[3]['S']
$ws.Cells.Item(11, 20).Value = $orig[3]['T']

# Row 12 <- original row 4
$ws.Cells.Item(12, 4).Value = $orig[4]['D']
$ws.Cells.Item(12, 11).Value = $orig[4]['K']
$ws.Cells.Item(12, 12).Value = $orig[4]['L']
$ws.Cells.Item(12, 13).Value = $orig[4]['M']
$ws.Cells.Item(12, 14).Value = $orig[4]['N']
$ws.Cells.Item(12, 15).Value = $orig[4]['O']
$ws.Cells.Item(12, 16).Value = $orig[4]['P']
$ws.Cells.Item(12, 17).Value = $orig[4]['Q']
$ws.Cells.Item(12, 18).Value = $orig[4]['R']
$ws.Cells.Item(12, 19).Value = $orig[4]['S']
$ws.Cells.Item(12, 20).Value = $orig[4]['T']

# Row 13 <- original row 2
$ws.Cells.Item(13, 4).Value = $orig[2]['D']
$ws.Cells.Item(13, 11).Value = $orig[2]['K']
$ws.Cells.Item(13, 12).Value = $orig[2]['L']
$ws.Cells.Item(13, 13).Value = $orig[2]['M']
$ws.Cells.Item(13, 14).Value = $orig[2]['N']
$ws.Cells.Item(13, 15).Value = $orig[2]['O']
$ws.Cells.Item(13, 16).Value = $orig[2]['P']
$ws.Cells.Item(13, 17).Value = $orig[2]['Q']
$ws.Cells.Item(13, 18).Value = $orig[2]['R']
$ws.Cells.Item(13, 19).Value = $orig[2]['S']
$ws.Cells.Item(13, 20).Value = $orig[2]['T']

# Row 14 <- original row 10
$ws.Cells.Item(14, 4).Value = $orig[10]['D']
$ws.Cells.Item(14, 11).Value = $orig[10]['K']
$ws.Cells.Item(14, 12).Value = $orig[10]['L']
$ws.Cells.Item(14, 13).Value = $orig[10]['M']
$ws.Cells.Item(14, 14).Value = $orig[10]['N']
$ws.Cells.Item(14, 15).Value = $orig[10]['O']
$ws.Cells.Item(14, 16).Value = $orig[10]['P']
$ws.Cells.Item(14, 17).Value = $orig[10]['Q']
$ws.Cells.Item(14, 18).Value = $orig[10]['R']
$ws.Cells.Item(14, 19).Value = $orig[10]['S']
$ws.Cells.Item(14, 20).Value = $orig[10]['T']

# Row 15 <- original row 11
$ws.Cells.Item(15, 4).Value = $orig[11]['D']
$ws.Cells.Item(15, 11).Value = $orig[11]['K']
$ws.Cells.Item(15, 12).Value = $orig[11]['L']
$ws.Cells.Item(15, 13).Value = $orig[11]['M']
$ws.Cells.Item(15, 14).Value = $orig[11]['N']
$ws.Cells.Item(15, 15).Value = $orig[11]['O']
$ws.Cells.Item(15, 16).Value = $orig[11]['P']
$ws.Cells.Item(15, 17).Value = $orig[11]['Q']
$ws.Cells.Item(15, 18).Value = $orig[11]['R']
$ws.Cells.Item(15, 19).Value = $orig[11]['S']
$ws.Cells.Item(15, 20).Value = $orig[11]['T']

# Row 16 <- original row 7
$ws.Cells.Item(16, 4).Value = $orig[7]['D']
$ws.Cells.Item(16, 11).Value = $orig[7]['K']
$ws.Cells.Item(16, 12).Value = $orig[7]['L']
$ws.Cells.Item(16, 13).Value = $orig[7]['M']
$ws.Cells.Item(16, 14).Value = $orig[7]['N']
$ws.Cells.Item(16, 15).Value = $orig[7]['O']
$ws.Cells.Item(16, 16).Value = $orig[7]['P']
$ws.Cells.Item(16, 17).Value = $orig[7]['Q']
$ws.Cells.Item(16, 18).Value = $orig[7]['R']
$ws.Cells.Item(16, 19).Value = $orig[7]['S']
$ws.Cells.Item(16, 20).Value = $orig[7]['T']

# Row 17 <- original row 16
$ws.Cells.Item(17, 4).Value = $orig[16]['D']
$ws.Cells.Item(17, 11).Value = $orig[16]['K']
$ws.Cells.Item(17, 12).Value = $orig[16]['L']
$ws.Cells.Item(17, 13).Value = $orig[16]['M']
$ws.Cells.Item(17, 14).Value = $orig[16]['N']
$ws.Cells.Item(17, 15).Value = $orig[16]['O']
$ws.Cells.Item(17, 16).Value = $orig[16]['P']
$ws.Cells.Item(17, 17).Value = $orig[16]['Q']
$ws.Cells.Item(17, 18).Value = $orig[16]['R']
$ws.Cells.Item(17, 19).Value = $orig[16]['S']
$ws.Cells.Item(17, 20).Value = $orig[16]['T']

# Row 18 <- original row 5
$ws.Cells.Item(18, 4).Value = $orig[5]['D']
$ws.Cells.Item(18, 11).Value = $orig[5]['K']
$ws.Cells.Item(18, 12).Value = $orig[5]['L']
$ws.Cells.Item(18, 13).Value = $orig[5]['M']
$ws.Cells.Item(18, 14).Value = $orig[5]['N']
$ws.Cells.Item(18, 15).Value = $orig[5]['O']
$ws.Cells.Item(18, 16).Value = $orig[5]['P']
$ws.Cells.Item(18, 17).Value = $orig[5]['Q']
$ws.Cells.Item(18, 18).Value = $orig[5]['R']
$ws.Cells.Item(18, 19).Value = $orig[5]['S']
$ws.Cells.Item(18, 20).Value = $orig[5]['T']

# Row 19 <- original row 6
$ws.Cells.Item(19, 4).Value = $orig[6]['D']
$ws.Cells.Item(19, 11).Value = $orig[6]['K']
$ws.Cells.Item(19, 12).Value = $orig[6]['L']
$ws.Cells.Item(19, 13).Value = $orig[6]['M']
$ws.Cells.Item(19, 14).Value = $orig[6]['N']
$ws.Cells.Item(19, 15).Value = $orig[6]['O']
$ws.Cells.Item(19, 16).Value = $orig[6]['P']
$ws.Cells.Item(19, 17).Value = $orig[6]['Q']
$ws.Cells.Item(19, 18).Value = $orig[6]['R']
$ws.Cells.Item(19, 19).Value = $orig[6]['S']
$ws.Cells.Item(19, 20).Value = $orig[6]['T']

